# Updated the test data as per latest Decision Table matrix(Sprint3)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transmittals_New")

# --- Row 2: remove CC (B2); Action-Level2 changes from SPInstall -> Submission
$ws.Range("B2").ClearContents()
$ws.Range("M2").Value = "Submission"

# --- Row 3: To becomes AutoTestAdmin (was AutoTestAdmin@@AutoTestUser);
#     clear the AttachDocuments..ReviewDocument block (G3:K3);
#     Action-Level2 changes from SPInstall -> Submission
$ws.Range("A3").Value = "AutoTestAdmin"
$ws.Range("G3:K3").ClearContents()
$ws.Range("M3").Value = "Submission"

# --- Row 4: add back the AttachDocuments..ReviewDocument block (G4:K4);
#     Action-Level3 changes from SPInstall -> Submission
$ws.Range("G4").Value = "Document Register"
$ws.Range("H4").Value = "Test 1 ta.docx"
$ws.Range("I4").Value = "Document Register"
$ws.Range("J4").Value = "Test 1 ta.docx"
$ws.Range("K4").Value = "BrowseDocument.docx"
$ws.Range("O4").Value = "Submission"

# --- Row 5: To becomes AutoTestAdmin (was AutoTestAdmin@@AutoTestUser);
#     ForwardTo becomes AutoTestUser, Action-Level3 becomes Submission
$ws.Range("A5").Value = "AutoTestAdmin"
$ws.Range("N5").Value = "AutoTestUser"
$ws.Range("O5").Value = "Submission"

# --- Row 6 (new row)
$ws.Range("A6").Value = "AutoTestAdmin@@AutoTestUser"
$ws.Range("C6").Value = "New Transmittal from Automation"
$ws.Range("D6").Value = "UnTick"
$ws.Range("E6").Value = "Correspondence"
$ws.Range("F6").Value = "Request for Information"
$ws.Range("G6").Value = "Document Register"
$ws.Range("H6").Value = "Test 1 ta.docx"
$ws.Range("I6").Value = "Document Register"
$ws.Range("J6").Value = "Test 1 ta.docx"
$ws.Range("K6").Value = "BrowseDocument.docx"
$ws.Range("L6").Value = "Message for New transmittal"
$ws.Range("M6").Value = "Submission"

# --- Row 7 (new row)
$ws.Range("A7").Value = "AutoTestAdmin@@AutoTestUser"
$ws.Range("C7").Value = "New Transmittal from Automation"
$ws.Range("D7").Value = "UnTick"
$ws.Range("E7").Value = "Correspondence"
$ws.Range("F7").Value = "Request for Information"
$ws.Range("L7").Value = "Message for New transmittal"
$ws.Range("M7").Value = "Submission"

# --- Sheet2: Action-Level2 on row 2 changes from SPInstall -> Submission
$ws2 = $wb.Worksheets.Item("Transmittals_New_ActionRequired")
$ws2.Range("M2").Value = "Submission"
